$wb = $excel.ActiveWorkbook

# --- constants sheet: delete "start_time" row (row 5) and "plot_end_time" row (old row 9) ---
$wsConstants = $wb.Worksheets.Item("constants")
$wsConstants.Rows.Item(9).Delete()
$wsConstants.Rows.Item(5).Delete()

# --- time_variants sheet: delete "freeze_times" row (row 14) ---
$wsTimeVariants = $wb.Worksheets.Item("time_variants")
$wsTimeVariants.Rows.Item(14).Delete()

# --- selection / active sheet changes ---
$wsConstants.Range("A7").Select()
$wsTimeVariants.Range("A7").Select()
$wsTimeVariants.Activate()
